$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.827.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4528"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3490"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.09"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07345"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.088"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.59"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.974"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.165"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.759.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.79"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06434"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.741"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.856.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.154"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.960.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.154"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.076"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09256"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.640"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.538"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.72"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06086"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02251"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2071"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.888"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6194"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.369"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.770"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.727"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5792"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.924"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.122"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06780"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.29"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.73%  "
